# Workbook was opened from https://d.docs.live.net/3b81fdb01c0fb0e6/ "wizualizacja dane.xlsx"
# Target sheet is "Arkusz8" (8th tab, rId8 / sheet8.xml) - the one that was tabSelected/active.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz8")
$ws.Activate()

# The percentage column (B2:B6) was reformatted from percent (0%/0.00%) to
# plain General numbers, and the underlying figures were re-entered as the
# "times 100" whole/decimal numbers instead of fractions.
$rng = $ws.Range("B2:B6")
$rng.NumberFormat = "General"

$ws.Range("B2").Value = 12
$ws.Range("B3").Value = 17
$ws.Range("B4").Value = 28
$ws.Range("B5").Value = 29.4
$ws.Range("B6").Value = 32.3

# Page setup was touched (portrait, paper size 9 = A4) which is what adds the
# <pageSetup/> element to the sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Final selection ended up on B6.
$ws.Range("B6").Select()
